# Generate Report for Handback
# Updates the handoff/handback timestamps recorded for the 49a60d07 file
# in the zh-cn and de-de localization sheets, and rolls the corresponding
# "Latest HO Xliff Generate Date" forward on the Overview sheet.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: row 2 (49a60d07 file) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-19 00:46:13"
$wsZhCn.Range("K2").Value = "2016-08-19 00:46:29"

# --- de-de sheet: row 2 (49a60d07 file) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-19 00:46:19"
$wsDeDe.Range("K2").Value = "2016-08-19 00:46:36"

# --- Overview sheet: refresh "Latest HO Xliff Generate Date" for the
#     49a60d07 file row (row 2) to the newest handback datetime above ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-19 00:46:19"
